$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.007322581950575113
$ws.Range("E2").Value = 0.1018237750977278
$ws.Range("G2").Value = 0.007695870008319616
$ws.Range("H2").Value = 0.06681700050830841
$ws.Range("I2").Value = 0.002055973745882511
$ws.Range("J2").Value = 0.01516876369714737
$ws.Range("K2").Value = 0.002399732824414968
$ws.Range("D3").Value = 0.001035409513860941
$ws.Range("E3").Value = 0.3582811197265983
$ws.Range("G3").Value = 0.02867878833785653
$ws.Range("H3").Value = 0.2299934513866901
$ws.Range("I3").Value = 0.01015108777210116
$ws.Range("J3").Value = 0.04752691462635994
$ws.Range("K3").Value = 0.01020820811390877
$ws.Range("D4").Value = 0.008836568333208561
$ws.Range("E4").Value = 0.1292784819379449
$ws.Range("G4").Value = 0.009095696732401848
$ws.Range("H4").Value = 0.08478125417605042
$ws.Range("I4").Value = 0.002608806826174259
$ws.Range("J4").Value = 0.01961160218343139
$ws.Range("K4").Value = 0.003228395711630583
$ws.Range("D5").Value = 0.001646819058805704
$ws.Range("E5").Value = 0.4029014273546636
$ws.Range("G5").Value = 0.03182772826403379
$ws.Range("H5").Value = 0.2593141803517938
$ws.Range("I5").Value = 0.01358586363494396
$ws.Range("J5").Value = 0.05136972526088357
$ws.Range("K5").Value = 0.01172052184119821
$ws.Range("E6").Value = 0.9538764897733927
$ws.Range("D7").Value = 0.005482353270053864
$ws.Range("E7").Value = 0.07444539573043585
$ws.Range("G7").Value = 0.004859731066972017
$ws.Range("H7").Value = 0.04732069931924343
$ws.Range("I7").Value = 0.001852010376751423
$ws.Range("J7").Value = 0.01368260616436601
$ws.Range("K7").Value = 0.001672299578785896
$ws.Range("D8").Value = 0.0008758548647165298
$ws.Range("E8").Value = 0.2877091160044074
$ws.Range("G8").Value = 0.02327139861881733
$ws.Range("H8").Value = 0.1859019221737981
$ws.Range("I8").Value = 0.007711600512266159
$ws.Range("J8").Value = 0.0370238465256989
$ws.Range("K8").Value = 0.00832445127889514
$ws.Range("D9").Value = 0.006412862800061703
$ws.Range("E9").Value = 0.09168151998892426
$ws.Range("G9").Value = 0.006185937207192183
$ws.Range("H9").Value = 0.05912051210179925
$ws.Range("I9").Value = 0.002011967822909355
$ws.Range("J9").Value = 0.01534664677456021
$ws.Range("K9").Value = 0.002233991865068674
$ws.Range("D10").Value = 0.001465102192014456
$ws.Range("E10").Value = 0.3535058093257248
$ws.Range("G10").Value = 0.02789804013445973
$ws.Range("H10").Value = 0.2279132469557226
$ws.Range("I10").Value = 0.01136909332126379
$ws.Range("J10").Value = 0.04438761109486222
$ws.Range("K10").Value = 0.01026897225528955
$ws.Range("E11").Value = 0.8631404815241694
$ws.Range("D12").Value = 0.003124778624624014
$ws.Range("E12").Value = 0.05461055412888527
$ws.Range("G12").Value = 0.003828324377536774
$ws.Range("H12").Value = 0.0354169849306345
$ws.Range("I12").Value = 0.001112679019570351
$ws.Range("J12").Value = 0.00840382743626833
$ws.Range("K12").Value = 0.001416357234120369
$ws.Range("D13").Value = 0.000872946809977293
$ws.Range("E13").Value = 0.2777701411396265
$ws.Range("G13").Value = 0.02208305709064007
$ws.Range("H13").Value = 0.1790438904426992
$ws.Range("I13").Value = 0.007510946597903967
$ws.Range("J13").Value = 0.0364212105050683
$ws.Range("K13").Value = 0.008111270610243082
$ws.Range("D14").Value = 0.002861146815121174
$ws.Range("E14").Value = 0.04968480579555035
$ws.Range("G14").Value = 0.003280565142631531
$ws.Range("H14").Value = 0.03195054223760962
$ws.Range("I14").Value = 0.001105108764022589
$ws.Range("J14").Value = 0.008586075156927109
$ws.Range("K14").Value = 0.001208995468914509
$ws.Range("D15").Value = 0.001343375537544489
$ws.Range("E15").Value = 0.2845779629424214
$ws.Range("G15").Value = 0.02225219598039985
$ws.Range("H15").Value = 0.1848116684705019
$ws.Range("I15").Value = 0.009134972468018532
$ws.Range("J15").Value = 0.03530557407066226
$ws.Range("K15").Value = 0.00816273083910346
$ws.Range("E16").Value = 0.9390419572591782
$ws.Range("D17").Value = 0.003192984033375978
$ws.Range("E17").Value = 0.05968657089397311
$ws.Range("G17").Value = 0.004223072435706854
$ws.Range("H17").Value = 0.03784104296937585
$ws.Range("I17").Value = 0.001122236251831055
$ws.Range("J17").Value = 0.01052592322230339
$ws.Range("K17").Value = 0.001527272164821625
$ws.Range("D18").Value = 0.0009333062916994095
$ws.Range("E18").Value = 0.3100540842860937
$ws.Range("G18").Value = 0.02490810491144657
$ws.Range("H18").Value = 0.2005890929140151
$ws.Range("I18").Value = 0.0086867930367589
$ws.Range("J18").Value = 0.03949838085100055
$ws.Range("K18").Value = 0.008871165569871664
$ws.Range("D19").Value = 0.002967767417430878
$ws.Range("E19").Value = 0.051557338796556
$ws.Range("G19").Value = 0.003638145979493856
$ws.Range("H19").Value = 0.03314039576798677
$ws.Range("I19").Value = 0.00109127489849925
$ws.Range("J19").Value = 0.008736977819353342
$ws.Range("K19").Value = 0.001243905164301395
$ws.Range("D20").Value = 0.001468064729124308
$ws.Range("E20").Value = 0.3119523008354008
$ws.Range("G20").Value = 0.02475325390696526
$ws.Range("H20").Value = 0.2021430507302284
$ws.Range("I20").Value = 0.009953280445188284
$ws.Range("J20").Value = 0.03843353502452374
$ws.Range("K20").Value = 0.008988030254840851
$ws.Range("E21").Value = 0.9581344211474061
$ws.Range("D22").Value = 0.002981758210808039
$ws.Range("E22").Value = 0.05135160824283957
$ws.Range("G22").Value = 0.003661919850856066
$ws.Range("H22").Value = 0.03342994069680572
$ws.Range("I22").Value = 0.001056331675499678
$ws.Range("J22").Value = 0.007929214742034674
$ws.Range("K22").Value = 0.001306978985667229
$ws.Range("D23").Value = 0.0008408799767494202
$ws.Range("E23").Value = 0.2793110520578921
$ws.Range("G23").Value = 0.02277254359796643
$ws.Range("H23").Value = 0.1806641155853868
$ws.Range("I23").Value = 0.007609287276864052
$ws.Range("J23").Value = 0.0350108789280057
$ws.Range("K23").Value = 0.008100660517811775
$ws.Range("D24").Value = 0.002797433640807867
$ws.Range("E24").Value = 0.04392475681379437
$ws.Range("G24").Value = 0.003027439117431641
$ws.Range("H24").Value = 0.02779966033995152
$ws.Range("I24").Value = 0.0009687785059213638
$ws.Range("J24").Value = 0.007670196704566479
$ws.Range("K24").Value = 0.001151410862803459
$ws.Range("D25").Value = 0.001083121635019779
$ws.Range("E25").Value = 0.2651075110770762
$ws.Range("G25").Value = 0.02058612275868654
$ws.Range("H25").Value = 0.1708462117239833
$ws.Range("I25").Value = 0.008810519706457853
$ws.Range("J25").Value = 0.0340973143465817
$ws.Range("K25").Value = 0.007779397536069155
$ws.Range("E26").Value = 1.1352184927091
$ws.Range("D27").Value = 0.004420891217887401
$ws.Range("E27").Value = 0.06838358081877231
$ws.Range("G27").Value = 0.004853783547878266
$ws.Range("H27").Value = 0.04416513368487358
$ws.Range("I27").Value = 0.001439846213907003
$ws.Range("J27").Value = 0.01114206705242395
$ws.Range("K27").Value = 0.001664528157562017
$ws.Range("D28").Value = 0.0009116794914007187
$ws.Range("E28").Value = 0.3026251026429236
$ws.Range("G28").Value = 0.02434277851134539
$ws.Range("H28").Value = 0.1952384945005178
$ws.Range("I28").Value = 0.008333943039178848
$ws.Range("J28").Value = 0.03909624628722667
$ws.Range("K28").Value = 0.008723151218146086
$ws.Range("D29").Value = 0.004775155801326036
$ws.Range("E29").Value = 0.07322538066655397
$ws.Range("G29").Value = 0.005045556835830212
$ws.Range("H29").Value = 0.04735847292467952
$ws.Range("I29").Value = 0.001557187363505363
$ws.Range("J29").Value = 0.01199029972776771
$ws.Range("K29").Value = 0.001813339814543724
$ws.Range("D30").Value = 0.001401296630501747
$ws.Range("E30").Value = 0.3236090023070574
$ws.Range("G30").Value = 0.02546346820890904
$ws.Range("H30").Value = 0.209005671646446
$ws.Range("I30").Value = 0.01057074591517448
$ws.Range("J30").Value = 0.0407187519595027
$ws.Range("K30").Value = 0.009383930545300245
$ws.Range("E31").Value = 0.9698823684826493
